$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.100.40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "'2.300.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'301.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").Value = "'99.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.55%  "
$ws.Range("D7").Value = "'0.504"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.514"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.57%  "
$ws.Range("D10").Value = "'36.16"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.45%  "
$ws.Range("D11").Value = "'0.0792"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").Value = "'18.67"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +10.86%  "
$ws.Range("E13").Value = "  +1.21%  "
$ws.Range("D14").Value = "'6.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.82%  "
$ws.Range("D15").Value = "'2.663.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "'2.269.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.09%  "
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("D18").Value = "'42.998.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("D19").Value = "'12.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.89%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "'0.0₃0906"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'6.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.56%  "
$ws.Range("D22").Value = "'67.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.17%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'236.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D24").Value = "'2.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +13.29%  "
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("D26").Value = "'2.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").Value = "'25.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.00%  "
$ws.Range("D28").Value = "'2.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("D29").Value = "'34.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.39%  "
$ws.Range("D30").Value = "'167.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("D31").Value = "'9.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").Value = "'5.04"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.54%  "
$ws.Range("D34").Value = "'17.78"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.60%  "
$ws.Range("D35").Value = "'4.67"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.39%  "
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("D37").Value = "'0.0692"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.46%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.101"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("D39").Value = "'1.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.48%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "'2.82"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("E41").Value = "  +0.60%  "
$ws.Range("D42").Value = "'2.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0291"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.86%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "'1.984.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.61%  "
$ws.Range("D45").Value = "'10.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.31%  "
$ws.Range("D46").Value = "'17.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("D47").Value = "'2.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.27%  "
$ws.Range("D48").Value = "'55.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.17%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "'2.530.62"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'1.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.04%  "
$ws.Range("D51").Value = "'70.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.10%  "
